$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("config")

# Add the new key/value pair for the business-exception error message
$ws.Range("A11").Value = "coin_be_error"
$ws.Range("B11").Value = "We couldn't load: {0}"

# Move the active selection, matching the saved cursor position
$ws.Range("A17").Select()
